$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before H ("Price" and everything after shifts right by one)
$ws.Columns("H").Insert()

# New header for the inserted column
$ws.Range("H1").Value = "Category3"

# Update Category1 / Category2 text (F2, G2) and set the new Category3 value (H2)
$ws.Range("F2").Value = "Kindle Books > Humor & Entertainment > Activities, Puzzles & Games > Crosswords"
$ws.Range("G2").Value = "Kindle Books > Large Print"
$ws.Range("H2").Value = "Kindle Books > Games > Crosswords"
